$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Date" column (B2:B7) with the new execution timestamps
# for the refreshed RAD test run (Existing Liability - Motor Fuel Tax).
$ws.Range("B2").Value = "Sun Oct 13 00:04:03 EDT 2024"
$ws.Range("B3").Value = "Sun Oct 13 00:04:17 EDT 2024"
$ws.Range("B4").Value = "Sun Oct 13 00:04:33 EDT 2024"
$ws.Range("B5").Value = "Sun Oct 13 00:04:49 EDT 2024"
$ws.Range("B6").Value = "Sun Oct 13 00:05:06 EDT 2024"
$ws.Range("B7").Value = "Sun Oct 13 00:05:21 EDT 2024"
